# Insert a new "Future Goals/FEATURES" slide into the deck.
#
# Per the target diff, a brand-new slide is inserted right after the
# "Another image" slide (position 7) and right before the "Explanation"
# slide (which was previously position 8). The new slide reuses the same
# "Title and Content" layout as the surrounding slides, has the title
# "Future Goals/FEATURES" and an otherwise empty content placeholder.
# All the other slides (and their content/order) are left untouched.

$p = $ppt.ActivePresentation

# "Title and Content" is the 2nd layout on the deck's single slide master
# (same layout used by the "Explanation" / "Explanation (cont.)" slides).
$layout = $p.SlideMaster.CustomLayouts.Item(2)

# Insert the new slide at index 8 (1-based), i.e. right after slide 7
# ("Another image") and before the old slide 8 ("Explanation"), which then
# gets pushed down to index 9.
$s = $p.Slides.AddSlide(8, $layout)

# Title placeholder text.
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Future Goals/FEATURES"

# Content placeholder is left empty (matches the target: a blank body).
